$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "License Information" (Heading2) paragraph -> plain paragraph with a
#    bold run reading "Aquifer Open Study Notes (Book Intros)"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(4)
if ($p1.Range.Text.TrimEnd([char]13, [char]7) -ne "License Information") {
    throw "Paragraph 4 text mismatch: " + $p1.Range.Text
}
$xml1 = '<w:p ' + $wns + '>' +
          '<w:pPr/>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr></w:r>' +
          '<w:r><w:rPr><w:b/><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr>' +
            '<w:t>Aquifer Open Study Notes (Book Intros)</w:t>' +
          '</w:r>' +
        '</w:p>'
[void]$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) License paragraph -> new wording, dropping the Tyndale House Publishers /
#    CC BY-SA hyperlinks entirely.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(5)
$xml2 = '<w:p ' + $wns + '>' +
          '<w:pPr/>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr></w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr>' +
            '<w:t xml:space="preserve">This work is an adaptation of </w:t>' +
          '</w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr>' +
            '<w:t>Tyndale Open Study Notes</w:t>' +
          '</w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr>' +
            '<w:t xml:space="preserve"> &#169; 2023 Tyndale House Publishers, licensed under the CC BY-SA 4.0 license. The adaptation, </w:t>' +
          '</w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr>' +
            '<w:t>Aquifer Open Study Notes</w:t>' +
          '</w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr>' +
            '<w:t>, was created by Mission Mutual and is also licensed under CC BY-SA 4.0.</w:t>' +
          '</w:r>' +
        '</w:p>'
[void]$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) "This PDF version is provided under the same license." -> new sentence
#    (rebuild the whole paragraph via InsertXML so the blank leading/trailing
#    runs around the text run survive - a plain Find/Replace or Range.Text
#    assignment collapses those empty sibling runs away.)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(6)
if ($p3.Range.Text.TrimEnd([char]13, [char]7) -ne "This PDF version is provided under the same license.") {
    throw "Paragraph 6 text mismatch: " + $p3.Range.Text
}
$xml3 = '<w:p ' + $wns + '>' +
          '<w:pPr/>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr></w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr>' +
            '<w:t>This resource has been adapted into multiple languages, including English, Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文).</w:t>' +
          '</w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr></w:r>' +
          '<w:r><w:rPr><w:lang w:val="zh_TW" w:bidi="zh_TW"/></w:rPr></w:r>' +
        '</w:p>'
[void]$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 4) Drop the two now-unused hyperlinks (Tyndale House Publishers / CC BY-SA)
#    in case any remain (defensive - paragraph 2 replacement above already
#    removed them, this is a no-op safety net).
# ---------------------------------------------------------------------------
for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
    $h = $d.Hyperlinks($i)
    if ($h.Address -eq "https://tyndaleopenresources.com/" -or `
        $h.Address -eq "https://creativecommons.org/licenses/by-sa/4.0/legalcode.en") {
        $h.Range.Delete()
    }
}

Write-Output "done"
